$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for handoff
#
# The "f34f748b-9cbd-4dfa-87d3-75f23b9cb40f.md" entry moves from "Handed
# back: in sync with en-US" to "Ready for handoff" (a new handoff cycle was
# kicked off for it), while the other two entries (ffff868d2531... and
# ffffff0ca5e934...) stay "Handed back" but shift up a row as the report is
# regenerated with the newly-handed-off file sorted last.
# ---------------------------------------------------------------------------

# ----------------------------- Overview sheet ------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value2 = "ffff868d2531-c139-461c-a010-804db51013f8.md"
$ws.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws.Range("C2").Value2 = "Handed back: in sync with en-US"

$ws.Range("A3").Value2 = "ffffff0ca5e934-cc8a-45aa-b2a3-d5f9be3bcdbd.md"
$ws.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws.Range("C3").Value2 = "Handed back: in sync with en-US"

$ws.Range("A4").Value2 = "f34f748b-9cbd-4dfa-87d3-75f23b9cb40f.md"
$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("C4").Value2 = "Ready for handoff"

# ------------------------------- zh-cn sheet -------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value2 = "ffff868d2531-c139-461c-a010-804db51013f8.md"
$ws.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws.Range("C2").Value2 = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.92d123faba7748170c7859b78b8858d0bf204f00.zh-cn.xlf"
$ws.Range("D2").Value2 = "2016-01-19 07:15:24"
$ws.Range("E2").Value2 = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.md"
$ws.Range("F2").Value2 = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.92d123faba7748170c7859b78b8858d0bf204f00.zh-cn.xlf"
$ws.Range("G2").Value2 = "2016-01-19 07:16:06"
$ws.Range("H2").Value2 = "Include"

$ws.Range("A3").Value2 = "ffffff0ca5e934-cc8a-45aa-b2a3-d5f9be3bcdbd.md"
$ws.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws.Range("C3").Value2 = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.92d123faba7748170c7859b78b8858d0bf204f00.zh-cn.xlf"
$ws.Range("D3").Value2 = "2016-01-19 07:15:24"
$ws.Range("E3").Value2 = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.md"
$ws.Range("F3").Value2 = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.92d123faba7748170c7859b78b8858d0bf204f00.zh-cn.xlf"
$ws.Range("G3").Value2 = "2016-01-19 07:16:06"
$ws.Range("H3").Value2 = "Include"

$ws.Range("A4").Value2 = "f34f748b-9cbd-4dfa-87d3-75f23b9cb40f.md"
$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("C4").Value2 = "f34f748b-9cbd-4dfa-87d3-75f23b9cb40f.595fa5a2a9496c6c9e5011c46e9d43eb47591dc9.zh-cn.xlf"
$ws.Range("D4").Value2 = "2016-01-19 07:19:04"
$ws.Range("E4").Value2 = "f34f748b-9cbd-4dfa-87d3-75f23b9cb40f.md"
$ws.Range("F4").Value2 = "f34f748b-9cbd-4dfa-87d3-75f23b9cb40f.595fa5a2a9496c6c9e5011c46e9d43eb47591dc9.zh-cn.xlf"
$ws.Range("G4").Value2 = "2016-01-19 07:18:05"
$ws.Range("H4").Value2 = "Include"

# ------------------------------- de-de sheet -------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value2 = "ffff868d2531-c139-461c-a010-804db51013f8.md"
$ws.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws.Range("C2").Value2 = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.92d123faba7748170c7859b78b8858d0bf204f00.de-de.xlf"
$ws.Range("D2").Value2 = "2016-01-19 07:15:35"
$ws.Range("E2").Value2 = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.md"
$ws.Range("F2").Value2 = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.92d123faba7748170c7859b78b8858d0bf204f00.de-de.xlf"
$ws.Range("G2").Value2 = "2016-01-19 07:16:23"
$ws.Range("H2").Value2 = "Include"

$ws.Range("A3").Value2 = "ffffff0ca5e934-cc8a-45aa-b2a3-d5f9be3bcdbd.md"
$ws.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws.Range("C3").Value2 = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.92d123faba7748170c7859b78b8858d0bf204f00.de-de.xlf"
$ws.Range("D3").Value2 = "2016-01-19 07:15:35"
$ws.Range("E3").Value2 = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.md"
$ws.Range("F3").Value2 = "6ef09566-1f3e-4049-a389-5ba7897ab1dd.92d123faba7748170c7859b78b8858d0bf204f00.de-de.xlf"
$ws.Range("G3").Value2 = "2016-01-19 07:16:23"
$ws.Range("H3").Value2 = "Include"

$ws.Range("A4").Value2 = "f34f748b-9cbd-4dfa-87d3-75f23b9cb40f.md"
$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("C4").Value2 = "f34f748b-9cbd-4dfa-87d3-75f23b9cb40f.595fa5a2a9496c6c9e5011c46e9d43eb47591dc9.de-de.xlf"
$ws.Range("D4").Value2 = "2016-01-19 07:19:14"
$ws.Range("E4").Value2 = "f34f748b-9cbd-4dfa-87d3-75f23b9cb40f.md"
$ws.Range("F4").Value2 = "f34f748b-9cbd-4dfa-87d3-75f23b9cb40f.595fa5a2a9496c6c9e5011c46e9d43eb47591dc9.de-de.xlf"
$ws.Range("G4").Value2 = "2016-01-19 07:18:22"
$ws.Range("H4").Value2 = "Include"
